# Update countries & provincias Spain
# - Re-sort a few country rows (swap the country name shown, since the
#   underlying row data was re-ranked by "Casos totales")
# - Refresh some numeric stats for the affected rows
# - Bump the "last updated" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country-name labels for rows whose rank changed -----------------

# Bulgaria (row 89) <-> Kirguistan (row 90)
$a89 = $ws.Range("A89").Value()
$a90 = $ws.Range("A90").Value()
$ws.Range("A89").Value = $a90
$ws.Range("A90").Value = $a89

# Santa Lucia (row 201) <-> Laos (row 202)
$a201 = $ws.Range("A201").Value()
$a202 = $ws.Range("A202").Value()
$ws.Range("A201").Value = $a202
$ws.Range("A202").Value = $a201

# Fiyi (row 203) <-> Dominica (row 204)
$a203 = $ws.Range("A203").Value()
$a204 = $ws.Range("A204").Value()
$ws.Range("A203").Value = $a204
$ws.Range("A204").Value = $a203

# Groenlandia (row 208) <-> Islas Malvinas (row 209)
$a208 = $ws.Range("A208").Value()
$a209 = $ws.Range("A209").Value()
$ws.Range("A208").Value = $a209
$ws.Range("A209").Value = $a208

# Seychelles (row 212) <-> Montserrat (row 213)
$a212 = $ws.Range("A212").Value()
$a213 = $ws.Range("A213").Value()
$ws.Range("A212").Value = $a213
$ws.Range("A213").Value = $a212

# --- Refresh numeric stats -------------------------------------------------

# Row 83 (El Salvador)
$ws.Range("D83").Value = 3330
$ws.Range("E83").Value = 2044
$ws.Range("G83").Value = 10
$ws.Range("H83").Value = 143

# Row 89 (now Kirguistan)
$ws.Range("C89").Value = 309
$ws.Range("D89").Value = 2212
$ws.Range("E89").Value = 2255
$ws.Range("G89").Value = 3
$ws.Range("H89").Value = 46

# Row 90 (now Bulgaria)
$ws.Range("B90").Value = 4513
$ws.Range("C90").Value = 0
$ws.Range("D90").Value = 2457
$ws.Range("E90").Value = 1841
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 215

# Row 136 (Georgia)
$ws.Range("B136").Value = 921
$ws.Range("C136").Value = 2
$ws.Range("D136").Value = 781
$ws.Range("E136").Value = 126

# Row 212 (now Seychelles)
$ws.Range("D212").Value = 10
$ws.Range("H212").Value = 1

# Row 213 (now Montserrat)
$ws.Range("D213").Value = 11
$ws.Range("H213").Value = 0

# --- Update "last refreshed" banner ----------------------------------------

$ws.Range("A1").Value = "Datos actualizados a 27 de Junio de 2020 a las 08:32"
